$d = $word.ActiveDocument

# --- Hunk 1: insert a new first paragraph with multiple runs before the
# existing "When thinking..." paragraph. ---
$firstOld = $d.Paragraphs.First
$firstOld.Range.InsertParagraphBefore()

$newPara1 = $d.Paragraphs.First
$newPara1.Range.Text = "I am a current second-year MD/PhD student at Vanderbilt University in Nashville, TN. The self-generated “Interest Map” above illustrates the many interconnected domains of my scientific and research interest domains. My overarching "

$insPoint = $d.Paragraphs.First.Range
$insPoint.Collapse(0)
$insPoint.MoveEnd(1, -1)
$insPoint.InsertAfter("scientific aims are to use")

$insPoint2 = $d.Paragraphs.First.Range
$insPoint2.Collapse(0)
$insPoint2.MoveEnd(1, -1)
$insPoint2.InsertAfter(" ")

$insPoint3 = $d.Paragraphs.First.Range
$insPoint3.Collapse(0)
$insPoint3.MoveEnd(1, -1)
$insPoint3.InsertAfter("computational and informatics-based approaches to genomic medicine")

$insPoint4 = $d.Paragraphs.First.Range
$insPoint4.Collapse(0)
$insPoint4.MoveEnd(1, -1)
$insPoint4.InsertAfter(". ")

$insPoint5 = $d.Paragraphs.First.Range
$insPoint5.Collapse(0)
$insPoint5.MoveEnd(1, -1)
$insPoint5.InsertAfter("Learn more about how I generated my interest map here. ")

Write-Host "Hunk1 done"

# --- Hunk 2: after "For each of the classes..." paragraph, insert:
#   "Medicine"
#   bullet: "Current second year " + "Medical Scientist Training Program..."
#   bullet: "Spent time in India..." + bookmark + "journal and in multiple local newspapers"
$forEachPara = $d.Paragraphs.Item(4)
$forEachPara.Range.InsertParagraphAfter()

$medicinePara = $d.Paragraphs.Item(5)
$medicinePara.Range.Text = "Medicine"

$medicinePara2 = $d.Paragraphs.Item(5)
$medicinePara2.Range.InsertParagraphAfter()

$bullet1 = $d.Paragraphs.Item(6)
$bullet1.Range.Text = "Current second year Medical Scientist Training Program (MSTP student) where I am exploring my clinical interests through the MD/PhD dual degree"
$bullet1.Range.ListFormat.ApplyBulletDefault()

$bullet1b = $d.Paragraphs.Item(6)
$bullet1b.Range.InsertParagraphAfter()

$bullet2 = $d.Paragraphs.Item(7)
$bullet2.Range.Text = "Spent time in India studying public health with collaborators at the All India Institute of Medical Science (AIIMS) to investigate barriers to health care for women. These findings were published in an Indian scientific journal and in multiple local newspapers"
$bullet2.Range.ListFormat.ApplyBulletDefault()

Write-Host "Hunk2 done"
Write-Host "Para count:"
Write-Host $d.Paragraphs.Count

Write-Host "--- paragraph listing ---"
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $t = $pp.Range.Text
    if ($t.Length -gt 60) { $t = $t.Substring(0,60) }
    Write-Host "$i : $t"
}
